# Add two new columns, "I0" (I) and "IF" (J), to the stats sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells (row 1) ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting (bold font, thin border, centered alignment) from the
# existing "IP" header cell (H1) onto the two new header cells so they match
# the look of the rest of the header row.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# --- Data cells (rows 2-62) ---
$iValues = @(7,9,1,9,7,6,6,6,8,9,8,8,8,8,6,8,10,7,8,7,9,8,9,8,8,10,8,8,8,8,6,7,8,7,7,6,7,7,6,9,7,8,8,7,7,7,8,7,7,8,8,8,8,7,8,7,9,8,9,8,4)
$jValues = @(8,9,1,9,8,7,6,7,9,9,8,8,8,8,6,8,10,7,8,7,9,8,9,8,8,11,8,8,8,8,8,7,9,7,7,7,7,7,7,9,8,8,8,7,7,7,8,7,7,8,8,8,8,8,8,7,9,8,9,8,4)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}
